$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('C2').Value = 'MSG: None

MSG: The decision was recorded as no movie selected.
'
$ws.Range('D2').Value = 'no_decision, '
$ws.Range('C3').Value = 'MSG: None

MSG: The function for no decision has been executed, indicating that no concrete choice of movie was made for Friday.
'
$ws.Range('D3').Value = 'no_decision, '
$ws.Range('C4').Value = 'MSG: None

MSG: The decision has been recorded as "no decision" for Friday''s movie.
'
$ws.Range('D4').Value = 'no_decision, '
$ws.Range('C5').Value = 'MSG: None

MSG: The decision has been recorded, indicating that no movie will be shown on Friday.
'
$ws.Range('D5').Value = 'no_decision, '
$ws.Range('C6').Value = 'MSG: None

MSG: No decision was made regarding which movie to show on Friday.
'
$ws.Range('D6').Value = 'no_decision, '
$ws.Range('C7').Value = 'MSG: None

MSG: The decision regarding the movie to be shown on Friday has resulted in no choice being made.
'
$ws.Range('D7').Value = 'no_decision, '
$ws.Range('C8').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Oppenheimer."
'
$ws.Range('D8').Value = 'Oppenheimer_was_selected, '
$ws.Range('C9').Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Range('D9').Value = 'no_decision, '
$ws.Range('C10').Value = 'MSG: None

MSG: The decision has been recorded, and no movie has been selected for showing on Friday.
'
$ws.Range('D10').Value = 'no_decision, '
$ws.Range('C11').Value = 'MSG: None

MSG: The decision has been recorded as "no decision" regarding which movie to show on Friday.
'
$ws.Range('D11').Value = 'no_decision, '
$ws.Range('C12').Value = 'MSG: None

MSG: The movie "Barbie" has been successfully selected for acquisition.
'
$ws.Range('D12').Value = 'Barbie_was_selected, '
$ws.Range('C13').Value = 'MSG: None

MSG: The decision-making process for Friday''s movie has concluded without a definitive selection.
'
$ws.Range('D13').Value = 'no_decision, '
$ws.Range('C14').Value = 'MSG: None

MSG: The decision process concluded without reaching a definitive choice of movie for Friday, so the outcome is recorded as no decision made.
'
$ws.Range('D14').Value = 'no_decision, '
$ws.Range('C15').Value = 'MSG: None

MSG: The decision has been recorded, and no movie will be shown on Friday.
'
$ws.Range('D15').Value = 'no_decision, '
$ws.Range('C16').Value = 'MSG: None

MSG: The decision has been recorded as "no decision."
'
$ws.Range('D16').Value = 'no_decision, '
$ws.Range('C17').Value = 'MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
'
$ws.Range('D17').Value = 'no_decision, '
$ws.Range('C18').Value = 'MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
'
$ws.Range('D18').Value = 'no_decision, '
$ws.Range('C19').Value = 'MSG: None

MSG: The decision to select a movie for Friday was not reached, so there will be no acquisition made at this time.
'
$ws.Range('D19').Value = 'no_decision, '
$ws.Range('C20').Value = 'MSG: None

MSG: The decision-making process concluded without an agreement on which movie to show on Friday, so no selection was made.
'
$ws.Range('D20').Value = 'no_decision, '
$ws.Range('C21').Value = 'MSG: None

MSG: The decision has been recorded: "Barbie" will be the movie shown on Friday.
'
$ws.Range('D21').Value = 'Barbie_was_selected, '
$ws.Range('C22').Value = 'MSG: None

MSG: The decision on which movie to show on Friday was not finalized, resulting in the conclusion that no decision was made.
'
$ws.Range('D22').Value = 'no_decision, '
$ws.Range('C23').Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie shown on Friday.
'
$ws.Range('D23').Value = 'Barbie_was_selected, '
$ws.Range('C24').Value = 'MSG: None

MSG: The rights for both movies have been successfully acquired.
'
$ws.Range('D24').Value = 'both_movies, '
$ws.Range('C25').Value = 'MSG: None

MSG: The committee did not reach a decision regarding which movie to show on Friday, so no movie rights will be acquired at this time.
'
$ws.Range('D25').Value = 'no_decision, '
$ws.Range('C26').Value = 'MSG: None

MSG: No decision was made regarding the movie selection for Friday.
'
$ws.Range('D26').Value = 'no_decision, '
$ws.Range('C27').Value = 'MSG: None

MSG: The rights to both movies have been successfully acquired.
'
$ws.Range('D27').Value = 'both_movies, '
$ws.Range('C28').Value = 'MSG: None

MSG: The decision has been recorded, and the rights for "Barbie" have been acquired for the showing on Friday.
'
$ws.Range('D28').Value = 'Barbie_was_selected, '
$ws.Range('C29').Value = 'MSG: None

MSG: The decision has been recorded as no decision being made regarding Friday''s movie.
'
$ws.Range('D29').Value = 'no_decision, '
$ws.Range('C30').Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Range('D30').Value = 'no_decision, '
$ws.Range('C31').Value = 'MSG: None

MSG: The decision has been recorded, confirming that there is no agreement on which movie to show on Friday.
'
$ws.Range('D31').Value = 'no_decision, '
$ws.Range('C32').Value = 'MSG: None

MSG: The decision concluded with no agreement on what movie to show on Friday.
'
$ws.Range('D32').Value = 'no_decision, '
$ws.Range('C33').Value = 'MSG: None

MSG: The decision has been made to acquire the rights to show "Barbie."
'
$ws.Range('D33').Value = 'Barbie_was_selected, '
$ws.Range('C34').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie" to be shown on Friday.
'
$ws.Range('D34').Value = 'Barbie_was_selected, '
$ws.Range('C35').Value = 'MSG: None

MSG: The committee has not reached a decision about which movie to show on Friday.
'
$ws.Range('D35').Value = 'no_decision, '
$ws.Range('C36').Value = 'MSG: None

MSG: The decision-making committee did not arrive at a clear decision regarding the movie to be shown on Friday.
'
$ws.Range('D36').Value = 'no_decision, '
$ws.Range('C37').Value = 'MSG: None

MSG: The decision has been recorded to acquire the rights for both movies.
'
$ws.Range('D37').Value = 'both_movies, '
$ws.Range('C38').Value = 'MSG: None

MSG: The decision has been recorded as no selection for the movie to be shown on Friday.
'
$ws.Range('D38').Value = 'no_decision, '
$ws.Range('C39').Value = 'MSG: None

MSG: The rights for both "Barbie" and "Oppenheimer" have been successfully acquired.
'
$ws.Range('D39').Value = 'both_movies, '
$ws.Range('C40').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range('D40').Value = 'Barbie_was_selected, '
$ws.Range('C41').Value = 'MSG: None

MSG: The committee has not reached a decision regarding the movie to be shown on Friday, so the no_decision function has been executed.
'
$ws.Range('D41').Value = 'no_decision, '
$ws.Range('C42').Value = 'MSG: None

MSG: The decision is that no movie was selected in this meeting.
'
$ws.Range('D42').Value = 'no_decision, '
$ws.Range('C43').Value = 'MSG: None

MSG: The decision about which movie to show on Friday resulted in no conclusion being reached.
'
$ws.Range('D43').Value = 'no_decision, '
$ws.Range('C44').Value = 'MSG: None

MSG: I have recorded the decision to acquire the rights for "Barbie" as the movie to be shown on Friday.
'
$ws.Range('D44').Value = 'Barbie_was_selected, '
$ws.Range('C45').Value = 'MSG: None

MSG: The rights to both movies will be acquired.
'
$ws.Range('D45').Value = 'both_movies, '
$ws.Range('C46').Value = 'MSG: None

MSG: The decision has been recorded, and the rights to "Barbie" will be acquired for Friday''s movie showing.
'
$ws.Range('D46').Value = 'Barbie_was_selected, '
$ws.Range('C47').Value = 'MSG: None

MSG: The rights for both movies have been acquired successfully.
'
$ws.Range('D47').Value = 'both_movies, '
$ws.Range('C48').Value = 'MSG: None

MSG: The decision has been recorded as no movie was definitively selected for Friday''s showing.
'
$ws.Range('D48').Value = 'no_decision, '
$ws.Range('C49').Value = 'MSG: None

MSG: The function for no decision has been called, indicating that no selection was made for Friday''s movie.
'
$ws.Range('D49').Value = 'no_decision, '
$ws.Range('C50').Value = 'MSG: None

MSG: The decision has been recorded as no movie being selected for Friday.
'
$ws.Range('D50').Value = 'no_decision, '
$ws.Range('C51').Value = 'MSG: None

MSG: Given the discussion lacks a definitive decision about the movie to be shown on Friday, the outcome is that no decision can be made.
'
$ws.Range('D51').Value = 'no_decision, '
$ws.Range('C52').Value = 'MSG: None

MSG: The decision has been recorded, and there is no movie choice for Friday.
'
$ws.Range('D52').Value = 'no_decision, '
$ws.Range('C53').Value = 'MSG: None

MSG: The rights to both movies have been acquired successfully.
'
$ws.Range('D53').Value = 'both_movies, '
$ws.Range('C54').Value = 'MSG: None

MSG: The decision-making process about the movie to show on Friday did not result in a consensus. Therefore, no decision has been made.
'
$ws.Range('D54').Value = 'no_decision, '
$ws.Range('C55').Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie" for the movie shown on Friday.
'
$ws.Range('D55').Value = 'Barbie_was_selected, '
$ws.Range('C56').Value = 'MSG: None

MSG: The decision has been made to show "Barbie" on Friday.
'
$ws.Range('D56').Value = 'Barbie_was_selected, '
$ws.Range('C57').Value = 'MSG: None

MSG: The decision has been recorded, and there was no selection made regarding the movies for Friday.
'
$ws.Range('D57').Value = 'no_decision, '
$ws.Range('C58').Value = 'MSG: None

MSG: The decision has been recorded as no decision regarding the movie for Friday.
'
$ws.Range('D58').Value = 'no_decision, '
$ws.Range('C59').Value = 'MSG: None

MSG: The decision has been successfully recorded to acquire the rights for "Oppenheimer."
'
$ws.Range('D59').Value = 'Oppenheimer_was_selected, '
$ws.Range('C60').Value = 'MSG: None

MSG: The decision to acquire the rights for "Oppenheimer" has been recorded successfully.
'
$ws.Range('D60').Value = 'Oppenheimer_was_selected, '
$ws.Range('C61').Value = 'MSG: None

MSG: I have decided to acquire the rights for "Barbie," as it was selected by the committee for Friday''s movie.
'
$ws.Range('D61').Value = 'Barbie_was_selected, '
$ws.Range('C62').Value = 'MSG: None

MSG: The decision has been made to acquire the rights to "Barbie."
'
$ws.Range('D62').Value = 'Barbie_was_selected, '
$ws.Range('C63').Value = 'MSG: None

MSG: The movie "Barbie" has been successfully selected for acquisition.
'
$ws.Range('D63').Value = 'Barbie_was_selected, '
$ws.Range('C64').Value = 'MSG: None

MSG: The decision regarding Friday''s movie has resulted in no determination being made.
'
$ws.Range('D64').Value = 'no_decision, '
$ws.Range('C65').Value = 'MSG: None

MSG: The decision about which movie to show on Friday was unresolved.
'
$ws.Range('D65').Value = 'no_decision, '
$ws.Range('C66').Value = 'MSG: None

MSG: The decision has been made to acquire the rights for "Barbie."
'
$ws.Range('D66').Value = 'Barbie_was_selected, '
$ws.Range('C67').Value = 'MSG: None

MSG: The decision about which movie to show on Friday was not reached.
'
$ws.Range('D67').Value = 'no_decision, '
$ws.Range('C68').Value = 'MSG: None

MSG: The decision has been made: no movie was selected for Friday.
'
$ws.Range('D68').Value = 'no_decision, '
